# Almighty Pushhhhh . . !!!
# Add GENDER, GOVT_ID and PERMISSION columns to the employee template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new columns -------------------------------------------------
# New "GENDER" column lands at D (pushes EMAIL..DEPARTMENT one to the right).
$ws.Columns("D").Insert()
# New "GOVT_ID" column lands at G, i.e. right before EMPCODE (now at column G).
$ws.Columns("G").Insert()

# --- Header row (row 1) -----------------------------------------------------
$ws.Range("D1").Value = "GENDER"
$ws.Range("G1").Value = "GOVT_ID"
$ws.Range("N1").Value = "PERMISSION"

# --- Data row (row 2) --------------------------------------------------------
$ws.Range("D2").Value = "MALE / FEMALE"
# Match the blue "hyperlink-like" font used by the e-mail cell (now E2).
$ws.Range("D2").Font.Color = $ws.Range("E2").Font.Color
$ws.Range("D2").Font.Name = $ws.Range("E2").Font.Name
$ws.Range("D2").Font.Size = $ws.Range("E2").Font.Size

$ws.Range("G2").Value = "Aadhar / PAN"
$ws.Range("N2").Value = "YES / NO"

# --- Hyperlink: the e-mail address moved from D2 to E2 ----------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:A@gmail.com", "", "", "A@gmail.com")

# --- View / selection state ---------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("G2").Select()
